$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.465.23"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.423.61"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.94"
$ws.Range("E5").Value = "  +0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.98"
$ws.Range("E6").Value = "  +3.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +5.81%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.423.00"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +2.27%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  +1.62%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("E12").Value = "  +1.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.017.47"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("E14").Value = "  +0.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.18"
$ws.Range("E15").Value = "  -1.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.367.60"
$ws.Range("E16").Value = "  +0.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.415.90"
$ws.Range("E18").Value = "  -0.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").Value = "  +1.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.10"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.60"
$ws.Range("E22").Value = "  -0.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.27"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24
$ws.Range("E24").Value = "  +8.29%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.536"
$ws.Range("E26").Value = "  +1.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  +2.85%  "

# Row 28
$ws.Range("E28").Value = "  +2.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.81"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("E31").Value = "  +0.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.44"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("E33").Value = "  +0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.05"
$ws.Range("E34").Value = "  +0.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.26"
$ws.Range("E35").Value = "  -1.88%  "

# Row 36
$ws.Range("E36").Value = "  +1.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.95"
$ws.Range("E37").Value = "  +1.90%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.871"
$ws.Range("E38").Value = "  -1.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.69"
$ws.Range("E39").Value = "  -4.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  +4.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  -1.34%  "

# Row 42
$ws.Range("E42").Value = "  +1.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.724.03"
$ws.Range("E43").Value = "  +0.99%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.28"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0688"
$ws.Range("E45").Value = "  +1.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.10"
$ws.Range("E46").Value = "  +5.35%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "337.20"
$ws.Range("E47").Value = "  +10.93%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.95"
$ws.Range("E48").Value = "  -0.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0288"
$ws.Range("E49").Value = "  +0.49%  "

# Row 50
$ws.Range("E50").Value = "  +4.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.10"
$ws.Range("E51").Value = "  +7.51%  "
